# Generate Report for Handback
# Update the "Correspond Handback DateTime" / "Latest HO Xliff Generate Date"
# timestamps to reflect a fresh report-generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-07-26 08:18:23"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2").Value = "2016-07-26 08:18:12"
$wsZhCn.Range("J2").Value = "2016-07-26 08:19:25"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("J2").Value = "2016-07-26 08:19:43"
